$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 68

# Force text storage (matches the rest of the sheet, which stores every
# value - dates, numbers, thousand-separated figures - as literal text)
# so Excel doesn't reinterpret "2025-05-08" as a date serial or "5,356"
# as a number.
$rng = $ws.Range("A$row`:J$row")
$rng.NumberFormat = "@"

$ws.Range("A$row").Value = "2025-05-08"
$ws.Range("B$row").Value = "38"
$ws.Range("C$row").Value = "37.28"
$ws.Range("D$row").Value = "1.03"
$ws.Range("E$row").Value = "0.27"
$ws.Range("F$row").Value = "0.09"
$ws.Range("G$row").Value = "5,356"
$ws.Range("H$row").Value = "8,019"
$ws.Range("I$row").Value = "8,069"
$ws.Range("J$row").Value = "7.2366"

# Drop the formatting we applied so the new row ends up with the same
# (default) cell style as every other row in the sheet.
$rng.ClearFormats()
